$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 199 (existing rows 199:287 shift down to 200:288)
$ws.Rows.Item(199).Insert()

# Populate the newly inserted row 199 with the new data record
$ws.Range("A199").Value = 7
$ws.Range("B199").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C199").Value = "Ñuble"
$ws.Range("D199").Value = 44523
$ws.Range("E199").Value = 16
$ws.Range("F199").Value = 100114014
$ws.Range("G199").Value = "Betarraga"
$ws.Range("H199").Value = "Sin especificar"
$ws.Range("I199").Value = "Primera"
$ws.Range("J199").Value = 200
$ws.Range("K199").Value = 700
$ws.Range("L199").Value = 800
$ws.Range("M199").Value = 750
$ws.Range("N199").Value = '$/paquete 5 unidades'
$ws.Range("O199").Value = "Región del Maule"
$ws.Range("P199").Value = 150
$ws.Range("Q199").Value = 5
$ws.Range("R199").Value = "Hortaliza"
